{"js": "// Replace each division-problem cell text in the practice table with its\n// new value, in document order. All old strings are unique except\n// \"83\u00f73=\" which appears twice (row 2 col 2 and row 4 col 2); those two\n// occurrences are disambiguated by search-result index (document order).\nconst replacements = [\n  [\"91\u00f77=\", \"15\u00f76=\"],\n  [\"52\u00f76=\", \"69\u00f73=\"],\n  [\"69\u00f76=\", \"56\u00f78=\"],\n  [\"98\u00f72=\", \"49\u00f77=\"],\n  [\"18\u00f77=\", \"69\u00f72=\"],\n  [\"39\u00f78=\", \"50\u00f75=\"],\n  [\"23\u00f75=\", \"56\u00f78=\"],\n  [\"84\u00f73=\", \"13\u00f72=\"],\n  [\"32\u00f75=\", \"64\u00f77=\"],\n  [\"79\u00f79=\", \"43\u00f75=\"],\n  [\"89\u00f76=\", \"76\u00f79=\"],\n  [\"19\u00f75=\", \"32\u00f76=\"],\n  [\"47\u00f72=\", \"92\u00f75=\"],\n  [\"73\u00f72=\", \"95\u00f77=\"],\n  [\"70\u00f76=\", \"25\u00f73=\"],\n  [\"82\u00f75=\", \"40\u00f74=\"],\n  [\"29\u00f76=\", \"52\u00f75=\"],\n  [\"55\u00f75=\", \"60\u00f76=\"],\n  [\"47\u00f77=\", \"63\u00f77=\"],\n  [\"66\u00f77=\", \"83\u00f74=\"],\n  [\"97\u00f78=\", \"39\u00f79=\"],\n  [\"90\u00f77=\", \"24\u00f79=\"],\n  [\"94\u00f77=\", \"46\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// \"83\u00f73=\" occurs twice (first -> \"72\u00f77=\", second -> \"72\u00f75=\"); replace in\n// document order using the search-result index.\nconst dup = body.search(\"83\u00f73=\", { matchCase: true, matchWholeWord: false });\ndup.load(\"items\");\nawait context.sync();\nconst dupNew = [\"72\u00f77=\", \"72\u00f75=\"];\nfor (let i = 0; i < dup.items.length; i++) {\n  dup.items[i].insertText(dupNew[i], Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the 5x5 grid of division problems (one-digit divisor practice\n# sheet) to the new set of problems, cell by cell, in document/reading\n# order (row by row, left to right). Only the data rows (1, 5, 9, 13, 17\n# of the 20-row table) contain problems; the other rows are blank answer\n# rows and are left untouched.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n  @(\"15\u00f76=\", \"69\u00f73=\", \"56\u00f78=\", \"49\u00f77=\", \"69\u00f72=\"),\n  @(\"72\u00f77=\", \"50\u00f75=\", \"56\u00f78=\", \"13\u00f72=\", \"64\u00f77=\"),\n  @(\"43\u00f75=\", \"76\u00f79=\", \"32\u00f76=\", \"92\u00f75=\", \"95\u00f77=\"),\n  @(\"25\u00f73=\", \"72\u00f75=\", \"40\u00f74=\", \"52\u00f75=\", \"60\u00f76=\"),\n  @(\"63\u00f77=\", \"83\u00f74=\", \"39\u00f79=\", \"24\u00f79=\", \"46\u00f74=\")\n)\n\n$dataRows = @(1, 5, 9, 13, 17)\n\nfor ($i = 0; $i -lt $dataRows.Length; $i++) {\n  $row = $dataRows[$i]\n  for ($col = 1; $col -le 5; $col++) {\n    $t.Cell($row, $col).Range.Text = $newValues[$i][$col - 1]\n  }\n}\n\nWrite-Output \"done\"\n"}
